# Apply the commit's changes:
# 1. Delete row 3 (the "Super Silk Taschentücher Würfelbox" entry), which shifts
#    every subsequent row up by one (old row 4 becomes row 3, ..., old row 36
#    disappears, leaving the sheet with rows 1-35 instead of 1-36).
# 2. Update the scrape timestamp (column O) for every data row (and the blank
#    row 2) from "2022-09-17 07:02:47" to "2022-09-17 20:59:34".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the obsolete product row; Excel automatically shifts subsequent rows
# up and updates the sheet's used range/dimension accordingly.
$ws.Rows(3).Delete()

# Refresh the recorded timestamp on every remaining data row (2 through 35).
$lastRow = $ws.UsedRange.Rows.Count
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 15).Value = "2022-09-17 20:59:34"
}
